# Applies the commit's change: adds a new column "param_E_pv3_solar"
# immediately before the existing "param_P_to_charging_station1" column
# (column U), shifting the two charging-station columns one place to the
# right, and updates the underlying data to the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at U; this shifts old U (param_P_to_charging_station1)
# to V and old V (param_P_to_charging_station2) to W, carrying their
# formatting and header styling along with them.
$ws.Columns("U:U").Insert()

# New header for the inserted column.
$ws.Range("U1").Value = "param_E_pv3_solar"

# New column data (constant 0.12, matching the other *_solar columns).
$ws.Range("U2").Value = 0.12
$ws.Range("U3").Value = 0.12
$ws.Range("U4").Value = 0.12
$ws.Range("U5").Value = 0.12
$ws.Range("U6").Value = 0.12
$ws.Range("U7").Value = 0.12
$ws.Range("U8").Value = 0.12
$ws.Range("U9").Value = 0.12
$ws.Range("U10").Value = 0.12
$ws.Range("U11").Value = 0.12
$ws.Range("U12").Value = 0.12
$ws.Range("U13").Value = 0.12
$ws.Range("U14").Value = 0.12
$ws.Range("U15").Value = 0.12
$ws.Range("U16").Value = 0.12
$ws.Range("U17").Value = 0.12

# Updated values for the (now shifted) charging-station columns V and W.
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0

$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0

$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0

$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 0

$ws.Range("V7").Value = 0
$ws.Range("W7").Value = 0

$ws.Range("V8").Value = 91.43994444444448
$ws.Range("W8").Value = 23.15744444444445

$ws.Range("V9").Value = 89.48500000000003
$ws.Range("W9").Value = 122.4135000000001

$ws.Range("V10").Value = 38.11500000000002
$ws.Range("W10").Value = 0

$ws.Range("V11").Value = 20.66350000000001
$ws.Range("W11").Value = 0

$ws.Range("V12").Value = 42.35
$ws.Range("W12").Value = 52.40552777777779

$ws.Range("V13").Value = 69.46500000000002
$ws.Range("W13").Value = 65.41700000000003

$ws.Range("V14").Value = 66.00000000000004
$ws.Range("W14").Value = 34.73616666666666

$ws.Range("V15").Value = 0
$ws.Range("W15").Value = 0

$ws.Range("V16").Value = 0
$ws.Range("W16").Value = 111.9891666666667

$ws.Range("V17").Value = 49.93404166666668
$ws.Range("W17").Value = 33
